$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 2
$ws.Range("G2").Value = 2.05
$ws.Range("H2").Value = 3.25
$ws.Range("I2").Value = 3.7
$ws.Range("J2").Value = 2.75
$ws.Range("K2").Value = 2.05
$ws.Range("Q2").Value = 2.25
$ws.Range("R2").Value = 1.62
$ws.Range("U2").Value = 1.95
$ws.Range("V2").Value = 1.8
$ws.Range("Y2").Value = 9
$ws.Range("AC2").Value = 8
$ws.Range("AD2").Value = 6.5
$ws.Range("AH2").Value = 9.5
$ws.Range("AP2").Value = 23
$ws.Range("BB2").Value = 251

# Row 4
$ws.Range("M4").Value = 1.06
$ws.Range("N4").Value = 10
$ws.Range("Q4").Value = 2.08
$ws.Range("R4").Value = 1.73

# Row 5
$ws.Range("O5").Value = 1.29
$ws.Range("P5").Value = 3.5

# Row 9
$ws.Range("G9").Value = 2.05
$ws.Range("H9").Value = 3.25
$ws.Range("I9").Value = 3.6
$ws.Range("J9").Value = 2.88
$ws.Range("K9").Value = 2
$ws.Range("L9").Value = 4.33
$ws.Range("M9").Value = 1.08
$ws.Range("N9").Value = 7.5
$ws.Range("Q9").Value = 2.3
$ws.Range("R9").Value = 1.6
$ws.Range("U9").Value = 2
$ws.Range("V9").Value = 1.73
$ws.Range("Y9").Value = 9.5
$ws.Range("Z9").Value = 19
$ws.Range("AA9").Value = 19
$ws.Range("AC9").Value = 7.5
$ws.Range("AF9").Value = 67
$ws.Range("AH9").Value = 9
$ws.Range("AI9").Value = 17
$ws.Range("AO9").Value = 12
$ws.Range("AP9").Value = 26
$ws.Range("AU9").Value = 9
$ws.Range("BB9").Value = 301

# Row 13
$ws.Range("M13").Value = 1.04
$ws.Range("N13").Value = 11.4

# Row 14
$ws.Range("G14").Value = 1.53
$ws.Range("H14").Value = 3.45
$ws.Range("I14").Value = 6.7
$ws.Range("J14").Value = 2.05
$ws.Range("K14").Value = 2.1
$ws.Range("L14").Value = 6.5
$ws.Range("N14").Value = 6.65
$ws.Range("O14").Value = 1.35
$ws.Range("P14").Value = 2.72
$ws.Range("U14").Value = 2.02
$ws.Range("V14").Value = 1.62
$ws.Range("W14").Value = 5.5
$ws.Range("X14").Value = 6.5
$ws.Range("Y14").Value = 8
$ws.Range("Z14").Value = 10.75
$ws.Range("AA14").Value = 13.5
$ws.Range("AB14").Value = 32
$ws.Range("AC14").Value = 7.9
$ws.Range("AD14").Value = 6.9
$ws.Range("AE14").Value = 19.5
$ws.Range("AH14").Value = 14
$ws.Range("AI14").Value = 40
$ws.Range("AM14").Value = 90
$ws.Range("AN14").Value = 3.25
$ws.Range("AO14").Value = 7.2
$ws.Range("AP14").Value = 17
$ws.Range("AR14").Value = 55
$ws.Range("AT14").Value = 2.47
$ws.Range("AU14").Value = 7.8
$ws.Range("AV14").Value = 80
$ws.Range("AX14").Value = 45
$ws.Range("AY14").Value = 45
$ws.Range("AZ14").Value = 350
$ws.Range("BA14").Value = 350
